# SAV-1059: add a second reference-data row to the "missing source" fixture.
# Row 3 deliberately omits sourceRecordId (column D) and category (column E)
# so the importer exercises the "missing source" validation path.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "InvoiceProduct-cat"
$ws.Range("B3").Value = "cat"
$ws.Range("C3").Value = "Cat"
# D3 (sourceRecordId) and E3 (category) are intentionally left blank.
$ws.Range("F3").Value = $true
$ws.Range("F3").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("G3").Value = "current"

$ws.Range("C10").Select() | Out-Null
